# The deck currently uses the "Integral" design (its colour scheme lives in
# the theme part backing the slide master).  The authored change swaps the
# presentation's theme palette back to the stock "Office Theme" colours
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), while the font scheme (Arial
# major/minor everywhere) and the fill/line/effect format scheme were
# already identical between the two themes, so only the twelve theme colours
# actually need to change.
#
# PowerPoint exposes the live theme colour scheme through
# Slide.ThemeColorScheme (it is shared by every slide/layout because they
# all hang off the single slide master), with slots ordered exactly like the
# MsoThemeColorSchemeIndex enum:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Office Theme palette (RGB() packs as R + G*256 + B*65536, matching the
# COM RGB() long values PowerPoint itself reports):
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$officeColors = @(
    0,         # 1  Dark1            000000
    16777215,  # 2  Light1           FFFFFF
    6968388,   # 3  Dark2            44546A
    15132391,  # 4  Light2           E7E6E6
    13998939,  # 5  Accent1          5B9BD5
    3243501,   # 6  Accent2          ED7D31
    10855845,  # 7  Accent3          A5A5A5
    49407,     # 8  Accent4          FFC000
    12874308,  # 9  Accent5          4472C4
    4697456,   # 10 Accent6          70AD47
    12673797,  # 11 Hyperlink        0563C1
    7491477    # 12 FollowedHyperlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
